# Generate Report for Handback
#
# The localization round-trip for 75b38656-db3a-4516-9551-dfd6146543e2.md
# has completed: the file's status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the "Latest Handback DateTime" moves
# forward to the new handback timestamp, and the stale "handback file is
# not the latest" error clears now that the handback is current.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns for the 75b38656... row
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus

# zh-cn detail sheet (table "zh_cn")
$ws_zhcn.Range("C3").Value = $newStatus
$ws_zhcn.Range("K3").Value = "2016-08-12 08:58:55"
$ws_zhcn.Range("P3").Value = ""
$ws_zhcn.Columns.Item(16).ColumnWidth = 12.9133

# de-de detail sheet (table "de_de")
$ws_dede.Range("C3").Value = $newStatus
$ws_dede.Range("K3").Value = "2016-08-12 08:59:09"
$ws_dede.Range("P3").Value = ""
$ws_dede.Columns.Item(16).ColumnWidth = 12.9133
